$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting bcc_emails/attachments/full_name/
# agent_code/leave_date columns one place to the right.
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "cc_emails"

# Match the hyperlink-style formatting used by the sibling email columns (A, C)
# for the new column's data rows, but leave them empty (no cc address yet).
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# Update the selection to mirror the author's cursor position after editing.
$ws.Range("D9").Select()
